$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Header change G1: END_PM -> End PM
$ws1.Range("G1").Value = "End PM"

# Replace numeric End PM values (G2:G19) with formatted inline strings,
# clearing the numeric right-aligned "Dialog" font style so the cells match
# the plain (unstyled) look of the Beg PM column.
$ws1.Range("G2").ClearFormats()
$ws1.Range("G2").Value = "End PM:   33.000"
$ws1.Range("G3").ClearFormats()
$ws1.Range("G3").Value = "End PM:   15.000"
$ws1.Range("G4").ClearFormats()
$ws1.Range("G4").Value = "End PM:   34.800"
$ws1.Range("G5").ClearFormats()
$ws1.Range("G5").Value = "End PM:   16.200"
$ws1.Range("G6").ClearFormats()
$ws1.Range("G6").Value = "End PM:   17.400"
$ws1.Range("G7").ClearFormats()
$ws1.Range("G7").Value = "End PM:   26.300"
$ws1.Range("G8").ClearFormats()
$ws1.Range("G8").Value = "End PM:   10.726"
$ws1.Range("G9").ClearFormats()
$ws1.Range("G9").Value = "End PM:   36.400"
$ws1.Range("G10").ClearFormats()
$ws1.Range("G10").Value = "End PM:   19.000"
$ws1.Range("G11").ClearFormats()
$ws1.Range("G11").Value = "End PM:   21.600"
$ws1.Range("G12").ClearFormats()
$ws1.Range("G12").Value = "End PM:   13.579"
$ws1.Range("G13").ClearFormats()
$ws1.Range("G13").Value = "End PM:   19.000"
$ws1.Range("G14").ClearFormats()
$ws1.Range("G14").Value = "End PM:   35.100"
$ws1.Range("G15").ClearFormats()
$ws1.Range("G15").Value = "End PM:   42.014"
$ws1.Range("G16").ClearFormats()
$ws1.Range("G16").Value = "End PM:   28.923"
$ws1.Range("G17").ClearFormats()
$ws1.Range("G17").Value = "End PM:   48.558"
$ws1.Range("G18").ClearFormats()
$ws1.Range("G18").Value = "End PM:   10.271"
$ws1.Range("G19").ClearFormats()
$ws1.Range("G19").Value = "End PM:   37.900"

# Update SQL text on sheet 2 (A2) to include End PM column
$ws2.Range("A2").Value = "select a.ea, a.treatment, a.county, a.route, a.year, ('Beg PM: ' || to_char(a.beg_pm, 990.999)) as `"Beg PM`", ('End PM: ' || to_char(a.end_pm, 990.999)) as `"End PM`", (a.end_pm-a.beg_pm) as length, a.budget_group from s1383currentl a `nwhere a.county = 'SM' `nunion  `nselect b.ea, b.treatment, b.county, b.route, b.year, ('Beg PM: ' || to_char(b.beg_pm, 990.999)) as `"Beg PM`",  ('End PM: ' || to_char(b.end_pm, 990.999)) as `"End PM`",  (b.end_pm-b.beg_pm) as length, b.budget_group from s1383historyl b `nwhere b.county = 'SM' `norder by year"

# Re-fit the row height back to standard; the engine auto-expands row
# height whenever a cell gains embedded line breaks.
$ws2.Rows.Item(2).AutoFit()
